# This script applies a weekly data refresh to the "Higo" (fig) price
# sheet: the fortnightly price records for each quality grade ("Primera"
# / "Segunda") are rotated forward to the next reporting date, carrying
# along their Volumen / Precio mínimo / Precio máximo / Precio promedio
# ponderado / Origen / Precio $/Kg values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the affected rows/columns (Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Origen, Precio $/Kg).
$updates = @(
    @{ Row = 2;  D = 44292; M = 25;  N = 16000; O = 16000; P = 16000; R = "Región Metropolitana"; S = 2286 },
    @{ Row = 3;  D = 44292; M = 30;  N = 15000; O = 15000; P = 15000; R = "Región Metropolitana"; S = 2143 },
    @{ Row = 4;  D = 44301; M = 100; N = 14000; O = 14000; P = 14000; R = "Región Metropolitana"; S = 2000 },
    @{ Row = 5;  D = 44301; M = 80;  N = 12000; O = 12000; P = 12000; R = "Región Metropolitana"; S = 1714 },
    @{ Row = 6;  D = 44320; M = 20;  N = 12000; O = 12000; P = 12000; R = "Región Metropolitana"; S = 1714 },
    @{ Row = 7;  D = 44320; M = 30;  N = 8000;  O = 8000;  P = 8000;  R = "Región Metropolitana"; S = 1143 },
    @{ Row = 8;  D = 44322; M = 45;  N = 12000; O = 12000; P = 12000; R = "Región Metropolitana"; S = 1714 },
    @{ Row = 9;  D = 44322; M = 80;  N = 8000;  O = 8000;  P = 8000;  R = "Región Metropolitana"; S = 1143 },
    @{ Row = 10; D = 44300; M = 100; N = 15000; O = 15000; P = 15000; R = "Región Metropolitana"; S = 2143 },
    @{ Row = 11; D = 44300; M = 80;  N = 12000; O = 12000; P = 12000; R = "Región Metropolitana"; S = 1714 },
    @{ Row = 12; D = 44299; M = 80;  N = 15000; O = 15000; P = 15000; R = "Provincia de Santiago"; S = 2143 },
    @{ Row = 13; D = 44299; M = 75;  N = 12000; O = 12000; P = 12000; R = "Provincia de Santiago"; S = 1714 },
    @{ Row = 14; D = 44302; M = 50;  N = 15000; O = 15000; P = 15000; R = "Región Metropolitana"; S = 2143 },
    @{ Row = 15; D = 44302; M = 30;  N = 12000; O = 12000; P = 12000; R = "Región Metropolitana"; S = 1714 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("D$r").Value = $u.D
    $ws.Range("M$r").Value = $u.M
    $ws.Range("N$r").Value = $u.N
    $ws.Range("O$r").Value = $u.O
    $ws.Range("P$r").Value = $u.P
    $ws.Range("R$r").Value = $u.R
    $ws.Range("S$r").Value = $u.S
}
